$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 180.7
$ws.Range("I6").Value = 145.22223
$ws.Range("K6").Value = 435.66669
$ws.Range("M6").Value = -323.66669
$ws.Range("H19").Value = 2170.9312
$ws.Range("I19").Value = 4254.846
$ws.Range("J19").Value = 477.75
$ws.Range("K19").Value = 4254.846
$ws.Range("L19").Value = 477.75
$ws.Range("M19").Value = -4079.846
$ws.Range("N19").Value = -827.75
$ws.Range("H33").Value = 880.4054
$ws.Range("I33").Value = 601.86957
$ws.Range("J33").Value = 1338
$ws.Range("K33").Value = 601.86957
$ws.Range("L33").Value = 1338
$ws.Range("M33").Value = -372.86957
$ws.Range("N33").Value = -1796
$ws.Range("H61").Value = 4496.5
$ws.Range("I61").Value = 94.5
$ws.Range("J61").Value = 8898.5
$ws.Range("K61").Value = 283.5
$ws.Range("L61").Value = 26695.5
$ws.Range("M61").Value = -111.5
$ws.Range("N61").Value = -27039.5
$ws.Range("H134").Value = 33437.5
$ws.Range("J134").Value = 33437.5
$ws.Range("L134").Value = 33437.5
$ws.Range("N134").Value = -43577.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 1186
$ws.Range("I6").Value = 776.6667
$ws.Range("J6").Value = 1800
$ws.Range("K6").Value = 776.6667
$ws.Range("L6").Value = 1800
$ws.Range("M6").Value = -603.6667
$ws.Range("N6").Value = -2146
$ws.Range("H32").Value = 9028.861999999999
$ws.Range("I32").Value = 5226.87
$ws.Range("J32").Value = 38304.2
$ws.Range("K32").Value = 5226.87
$ws.Range("L32").Value = 38304.2
$ws.Range("M32").Value = -4939.87
$ws.Range("N32").Value = -38878.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 393
$ws.Range("I8").Value = 393
$ws.Range("K8").Value = 393
$ws.Range("M8").Value = -253
$ws.Range("H86").Value = 4376.8823
$ws.Range("I86").Value = 4128.1055
$ws.Range("J86").Value = 4692
$ws.Range("K86").Value = 4128.1055
$ws.Range("L86").Value = 4692
$ws.Range("M86").Value = -3005.1055
$ws.Range("N86").Value = -6938
$ws.Range("H89").Value = 4376.8823
$ws.Range("I89").Value = 4128.1055
$ws.Range("J89").Value = 4692
$ws.Range("K89").Value = 20640.5275
$ws.Range("L89").Value = 23460
$ws.Range("M89").Value = -15024.5275
$ws.Range("N89").Value = -34692
$ws.Range("H98").Value = 30000
$ws.Range("J98").Value = 30000
$ws.Range("L98").Value = 30000
$ws.Range("N98").Value = -35990

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 3750
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 3750
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 3750
$ws.Range("M13").Value = $null
$ws.Range("N13").Value = -4028
$ws.Range("H58").Value = 2117.9312
$ws.Range("I58").Value = 1556
$ws.Range("J58").Value = 3366.6667
$ws.Range("K58").Value = 1556
$ws.Range("L58").Value = 3366.6667
$ws.Range("M58").Value = -1353
$ws.Range("N58").Value = -3772.6667
$ws.Range("H99").Value = 50370.145
$ws.Range("J99").Value = 2987.7144
$ws.Range("L99").Value = 2987.7144
$ws.Range("N99").Value = -5983.7144
$ws.Range("H122").Value = 1430.9535
$ws.Range("I122").Value = 1567.7693
$ws.Range("J122").Value = 1221.7059
$ws.Range("K122").Value = 4703.3079
$ws.Range("L122").Value = 3665.1177
$ws.Range("M122").Value = -2253.3079
$ws.Range("N122").Value = -8565.117699999999
$ws.Range("H126").Value = 50370.145
$ws.Range("J126").Value = 2987.7144
$ws.Range("L126").Value = 8963.143199999999
$ws.Range("N126").Value = -13903.1432
$ws.Range("H136").Value = 2117.9312
$ws.Range("I136").Value = 1556
$ws.Range("J136").Value = 3366.6667
$ws.Range("K136").Value = 4668
$ws.Range("L136").Value = 10100.0001
$ws.Range("M136").Value = -2118
$ws.Range("N136").Value = -15200.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 90909690
$ws.Range("I6").Value = 125000060
$ws.Range("J6").Value = 2033.3334
$ws.Range("K6").Value = 375000180
$ws.Range("L6").Value = 6100.0002
$ws.Range("M6").Value = -375000067
$ws.Range("N6").Value = -6326.0002
$ws.Range("H11").Value = 338
$ws.Range("J11").Value = 575
$ws.Range("L11").Value = 1725
$ws.Range("N11").Value = -2005
$ws.Range("H25").Value = 345
$ws.Range("I25").Value = 345
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 1035
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -866
$ws.Range("N25").Value = $null
$ws.Range("H30").Value = 345
$ws.Range("I30").Value = 345
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 1035
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -933
$ws.Range("N30").Value = $null
$ws.Range("H81").Value = 5940
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 5940
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 17820
$ws.Range("M81").Value = $null
$ws.Range("N81").Value = -20066
$ws.Range("H84").Value = 5940
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 5940
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 53460
$ws.Range("M84").Value = $null
$ws.Range("N84").Value = -64692
$ws.Range("H116").Value = 3480.2856
$ws.Range("I116").Value = 776.6667
$ws.Range("J116").Value = 5508
$ws.Range("K116").Value = 2330.0001
$ws.Range("L116").Value = 16524
$ws.Range("M116").Value = 1111.9999
$ws.Range("N116").Value = -23408

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6746.3657
$ws.Range("I70").Value = 3842.4546
$ws.Range("J70").Value = 18725
$ws.Range("K70").Value = 3842.4546
$ws.Range("L70").Value = 18725
$ws.Range("M70").Value = -3572.4546
$ws.Range("N70").Value = -19265
$ws.Range("H73").Value = 6746.3657
$ws.Range("I73").Value = 3842.4546
$ws.Range("J73").Value = 18725
$ws.Range("K73").Value = 3842.4546
$ws.Range("L73").Value = 18725
$ws.Range("M73").Value = -2906.4546
$ws.Range("N73").Value = -20597
